$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 (last existing header cell) onto the
# new H1 header cell, then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column (rows 2-8) with 0, matching the plain
# (unstyled) numeric cells used elsewhere in the data rows.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
